# Jun's Oct 9th updates to MN data
# Set the BIEfIE control-lever boolean value (cell B2 on the "BIEfIE" sheet) to 0,
# i.e. exclude emissions from imported electricity.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("BIEfIE")
$ws.Range("B2").Value = 0
